# This script rewrites the "hasInterval" property row (row 4) into a
# "hasInteger" property row, matching the target diff:
#   A4: hasInterval        -> hasInteger
#   B4: Time interval      -> has Integer
#   C4: Zeitintervall      -> Zahl
#   G4: Time interval      -> integer
#   H4: Zeitintervall      -> zahl
#   L4: hasSequenceBounds  -> hasValue
#   M4: IntervalValue      -> IntValue
#   N4: SimpleText         -> Spinbox
#   O4: maxlength: 5, rows: 10 -> max: 10, min: 5, rows: 10

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "hasInteger"
$ws.Range("B4").Value = "has Integer"
$ws.Range("C4").Value = "Zahl"
$ws.Range("G4").Value = "integer"
$ws.Range("H4").Value = "zahl"
$ws.Range("L4").Value = "hasValue"
$ws.Range("M4").Value = "IntValue"
$ws.Range("N4").Value = "Spinbox"
$ws.Range("O4").Value = "max: 10, min: 5, rows: 10"
